$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold numeric-looking values stored as TEXT (shared strings),
# e.g. " 1.437" with a leading space for alignment. A plain .Value
# assignment gets auto-coerced to a real number by Excel, so the cells
# that don't already carry a Text number format (column F already does)
# need one applied first; we then restore the original "Normal" style so
# the cell's style index is unaffected by the temporary format change.

# Row 2 (intrcpt): Chi2 (D2) and p value (F2)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = " 1.436"
$ws.Range("D2").Style = "Normal"
$ws.Range("F2").Value = "0.1511"

# Row 3 (GenLength_y_IUCN.y): Chi2 (D3) and p value (F3)
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = " 0.069"
$ws.Range("D3").Style = "Normal"
$ws.Range("F3").Value = "0.7925"

# Row 4 (Pvalue): Estimate (B4), Chi2 (D4), p value (F4)
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "-0.057"
$ws.Range("B4").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "-1.425"
$ws.Range("D4").Style = "Normal"

$ws.Range("F4").Value = "0.1542"
